$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts existing rows 13:37 down to 14:38)
$ws.Rows.Item(13).Insert()

# Populate the new row 13 with the weekly record
$ws.Cells.Item(13, 1).Value = 7
$ws.Cells.Item(13, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(13, 3).Value = "Ñuble"
$ws.Cells.Item(13, 4).Value = 44544
$ws.Cells.Item(13, 5).Value = 16
$ws.Cells.Item(13, 6).Value = 100112022
$ws.Cells.Item(13, 7).Value = "Arveja Verde"
$ws.Cells.Item(13, 8).Value = "Sin especificar"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 100
$ws.Cells.Item(13, 11).Value = 13000
$ws.Cells.Item(13, 12).Value = 14000
$ws.Cells.Item(13, 13).Value = 13500
$ws.Cells.Item(13, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(13, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(13, 16).Value = 540
$ws.Cells.Item(13, 17).Value = 25
$ws.Cells.Item(13, 18).Value = "Hortaliza"
